# Replace the "Meus Objetivos" body paragraph (and the blank paragraph that
# follows it) with the new merged wording + indentation, moving the
# "_GoBack" bookmark into the middle of the new text and dropping the old
# "_gjdgxs" bookmark.
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Trabalhar na área de desenvolvimento de software*") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Meus Objetivos' paragraph to replace."
}

$targetIndex = $target.Index
$nextPara = $d.Paragraphs.Item($targetIndex + 1)

# Consume both the text paragraph and the following empty paragraph so the
# two collapse into a single paragraph, matching the diff.
$combined = $d.Range($target.Range.Start, $nextPara.Range.End)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Sempre estar em constante aprendizado com a área de TI, atuando ativamente para o meu crescimento profissional e</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>dos meus parceiros (as), com foco sempre no Desenvolvimento S</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>eguro de aplicações WEB e CyberSecurity</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$combined.InsertXML($xmlFragment)
